$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "Please indicate who provides instruction to students while at outdoor school:`n (Instruction: structured or unstructured activities with explicit learning objectives) - 1 - "

$ws.Range("BS1").Value = $prefix + "adult volunteers"
$ws.Range("BT1").Value = $prefix + "parents"
$ws.Range("BU1").Value = $prefix + "college volunteers"
$ws.Range("BV1").Value = $prefix + "school teachers"
$ws.Range("BW1").Value = $prefix + "school administration"
$ws.Range("BX1").Value = $prefix + "high school volunteers"
$ws.Range("BY1").Value = $prefix + "trained staff"
$ws.Range("BZ1").Value = $prefix + "natural resource professionals"
$ws.Range("CA1").Value = $prefix + "other"
$ws.Range("CB1").Value = $prefix + "other - Text"
$ws.Range("CC1").Value = $prefix + "other"
$ws.Range("CD1").Value = $prefix + "other - Text"
$ws.Range("CE1").Value = $prefix + "other"
$ws.Range("CF1").Value = $prefix + "other - Text"
